# SSHConfig.xlsx: point the hostname var at the new KTH map target box
# (adeye05u -> adeye06u) and leave the cursor where the edit finished.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "hostname" row's value cell (A7 label / B7 value) -> bump the host suffix
$ws.Range("B7").Value = "adeye06u"

# Shrink the sheet-tab-area/horizontal-scrollbar split in the window chrome.
$excel.ActiveWindow.TabRatio = 0.5

# Leave selection where the author's cursor ended up after the edit.
$ws.Range("C9").Select()
